# Update table header: "P. BOLSA"/"P. MUESTRA" become just "BOLSA"/"MUESTRA".
# Shared-string table order matters (MUESTRA gets appended before BOLSA),
# so write D1 before C1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "MUESTRA"
$ws.Range("C1").Value = "BOLSA"

# Header row shrinks back to the standard row height used elsewhere in the sheet.
$ws.Rows.Item(1).RowHeight = 15.75

# Leave the selection on C1, matching the saved workbook state.
$ws.Range("C1").Select()
